# Applies scheduled-runner market price refresh to the Leve profit tables.
# Source data: FFXIV marketboard snapshot update (chore: update Sheets via scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3318.889
$ws.Range("J32").Value = 3576.2
$ws.Range("L32").Value = 3576.2
$ws.Range("N32").Value = -4228.2
$ws.Range("H53").Value = 294.33334
$ws.Range("I53").Value = 228.85715
$ws.Range("J53").Value = 523.5
$ws.Range("K53").Value = 228.85715
$ws.Range("L53").Value = 523.5
$ws.Range("M53").Value = 408.14285
$ws.Range("N53").Value = -1797.5
$ws.Range("H62").Value = 2999.6667
$ws.Range("I62").Value = 2999.6667
$ws.Range("K62").Value = 2999.6667
$ws.Range("M62").Value = -2375.6667
$ws.Range("H65").Value = 2999.6667
$ws.Range("I65").Value = 2999.6667
$ws.Range("K65").Value = 14998.3335
$ws.Range("M65").Value = -11878.3335
$ws.Range("H92").Value = 642.0769
$ws.Range("I92").Value = 601
$ws.Range("K92").Value = 601
$ws.Range("M92").Value = 647
$ws.Range("H98").Value = 1791
$ws.Range("I98").Value = 1831.4286
$ws.Range("J98").Value = 1649.5
$ws.Range("K98").Value = 1831.4286
$ws.Range("L98").Value = 1649.5
$ws.Range("M98").Value = -333.4286
$ws.Range("N98").Value = -4645.5
$ws.Range("H113").Value = 2865.125
$ws.Range("I113").Value = 1327.5
$ws.Range("J113").Value = 3377.6667
$ws.Range("K113").Value = 1327.5
$ws.Range("L113").Value = 3377.6667
$ws.Range("M113").Value = 1926.5
$ws.Range("N113").Value = -9885.6667
$ws.Range("H122").Value = 1791
$ws.Range("I122").Value = 1831.4286
$ws.Range("J122").Value = 1649.5
$ws.Range("K122").Value = 5494.2858
$ws.Range("L122").Value = 4948.5
$ws.Range("M122").Value = -3044.2858
$ws.Range("N122").Value = -9848.5
$ws.Range("H135").Value = 1398.85
$ws.Range("I135").Value = 1546.8
$ws.Range("J135").Value = 955
$ws.Range("K135").Value = 13921.2
$ws.Range("L135").Value = 8595
$ws.Range("M135").Value = -11386.2
$ws.Range("N135").Value = -13665
$ws.Range("H138").Value = 4746.8623
$ws.Range("J138").Value = 5692.625
$ws.Range("L138").Value = 17077.875
$ws.Range("N138").Value = -27357.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1665.2106
$ws.Range("I2").Value = 1212.3
$ws.Range("J2").Value = 2168.4443
$ws.Range("K2").Value = 1212.3
$ws.Range("L2").Value = 2168.4443
$ws.Range("M2").Value = -1099.3
$ws.Range("N2").Value = -2394.4443
$ws.Range("H74").Value = 1535.8
$ws.Range("I74").Value = 1419.75
$ws.Range("K74").Value = 1419.75
$ws.Range("M74").Value = -545.75
$ws.Range("H77").Value = 1535.8
$ws.Range("I77").Value = 1419.75
$ws.Range("K77").Value = 7098.75
$ws.Range("M77").Value = -2730.75
$ws.Range("H116").Value = 1665.2106
$ws.Range("I116").Value = 1212.3
$ws.Range("J116").Value = 2168.4443
$ws.Range("K116").Value = 1212.3
$ws.Range("L116").Value = 2168.4443
$ws.Range("M116").Value = 1081.7
$ws.Range("N116").Value = -6756.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1665.2106
$ws.Range("I3").Value = 1212.3
$ws.Range("J3").Value = 2168.4443
$ws.Range("K3").Value = 1212.3
$ws.Range("L3").Value = 2168.4443
$ws.Range("M3").Value = -1098.3
$ws.Range("N3").Value = -2396.4443
$ws.Range("H20").Value = 7427.5
$ws.Range("I20").Value = 7091.2
$ws.Range("J20").Value = 7988
$ws.Range("K20").Value = 7091.2
$ws.Range("L20").Value = 7988
$ws.Range("M20").Value = -6844.2
$ws.Range("N20").Value = -8482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2121.1853
$ws.Range("I31").Value = 1521
$ws.Range("K31").Value = 1521
$ws.Range("M31").Value = -1226
$ws.Range("H34").Value = 2121.1853
$ws.Range("I34").Value = 1521
$ws.Range("K34").Value = 1521
$ws.Range("M34").Value = -1319
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H86").Value = 5385
$ws.Range("I86").Value = 5300
$ws.Range("J86").Value = 5555
$ws.Range("K86").Value = 5300
$ws.Range("L86").Value = 5555
$ws.Range("M86").Value = -4177
$ws.Range("N86").Value = -7801
$ws.Range("H89").Value = 5385
$ws.Range("I89").Value = 5300
$ws.Range("J89").Value = 5555
$ws.Range("K89").Value = 26500
$ws.Range("L89").Value = 27775
$ws.Range("M89").Value = -20884
$ws.Range("N89").Value = -39007
$ws.Range("H107").Value = 1081.4
$ws.Range("I107").Value = 1051.75
$ws.Range("K107").Value = 1051.75
$ws.Range("M107").Value = 868.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2295.9333
$ws.Range("I68").Value = 1474.3334
$ws.Range("J68").Value = 2501.3333
$ws.Range("K68").Value = 4423.0002
$ws.Range("L68").Value = 7503.999899999999
$ws.Range("M68").Value = -3612.0002
$ws.Range("N68").Value = -9125.999899999999
$ws.Range("H71").Value = 2295.9333
$ws.Range("I71").Value = 1474.3334
$ws.Range("J71").Value = 2501.3333
$ws.Range("K71").Value = 13269.0006
$ws.Range("L71").Value = 22511.9997
$ws.Range("M71").Value = -9213.000599999999
$ws.Range("N71").Value = -30623.9997
$ws.Range("H107").Value = 1798.0834
$ws.Range("I107").Value = 1541.4286
$ws.Range("K107").Value = 4624.2858
$ws.Range("M107").Value = -2704.2858
$ws.Range("H113").Value = 1009.7647
$ws.Range("I113").Value = 842.5
$ws.Range("J113").Value = 1032.0667
$ws.Range("K113").Value = 2527.5
$ws.Range("L113").Value = 3096.2001
$ws.Range("M113").Value = -357.5
$ws.Range("N113").Value = -7436.2001
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 70188.07000000001
$ws.Range("I139").Value = 94029.55
$ws.Range("K139").Value = 282088.65
$ws.Range("M139").Value = -276948.65

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33959.8
$ws.Range("I80").Value = 13333
$ws.Range("K80").Value = 13333
$ws.Range("M80").Value = -12335
$ws.Range("H83").Value = 33959.8
$ws.Range("I83").Value = 13333
$ws.Range("K83").Value = 66665
$ws.Range("M83").Value = -61673
$ws.Range("H102").Value = 1911.1428
$ws.Range("I102").Value = 1884.7778
$ws.Range("K102").Value = 1884.7778
$ws.Range("M102").Value = -262.7778000000001
$ws.Range("H113").Value = 9400
$ws.Range("I113").Value = 6600
$ws.Range("K113").Value = 6600
$ws.Range("M113").Value = -4430
$ws.Range("H122").Value = 4009.3333
$ws.Range("I122").Value = 4132.6
$ws.Range("K122").Value = 12397.8
$ws.Range("M122").Value = -9947.800000000001
$ws.Range("H126").Value = 6392.7144
$ws.Range("J126").Value = 5708
$ws.Range("L126").Value = 17124
$ws.Range("N126").Value = -22064
$ws.Range("H132").Value = 8000
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1115.75
$ws.Range("I82").Value = 1041.5714
$ws.Range("K82").Value = 1041.5714
$ws.Range("M82").Value = -680.5714
$ws.Range("H85").Value = 1115.75
$ws.Range("I85").Value = 1041.5714
$ws.Range("K85").Value = 1041.5714
$ws.Range("M85").Value = 206.4286
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10000002
$ws.Range("J81").Value = 10000002
$ws.Range("L81").Value = 20000004
$ws.Range("N81").Value = -20002126
$ws.Range("H84").Value = 10000002
$ws.Range("J84").Value = 10000002
$ws.Range("L84").Value = 100000020
$ws.Range("N84").Value = -100010628
$ws.Range("H113").Value = 279.875
$ws.Range("I113").Value = 290
$ws.Range("K113").Value = 870
$ws.Range("M113").Value = 1300
$ws.Range("H122").Value = 7006.5713
$ws.Range("I122").Value = 5825.3335
$ws.Range("K122").Value = 17476.0005
$ws.Range("M122").Value = -15026.0005
$ws.Range("H132").Value = 564.9
$ws.Range("I132").Value = 483.125
$ws.Range("J132").Value = 892
$ws.Range("K132").Value = 1449.375
$ws.Range("L132").Value = 2676
$ws.Range("M132").Value = 1080.625
$ws.Range("N132").Value = -7736

